# Update "想去人数" (F column) figures for several events on the
# "展览" sheet and the corresponding rows on the "全部类型" sheet.
# These sheets list the same events (全部类型 is a superset that also
# includes the "演出" sheet entry), so the row numbers are offset by 1
# starting from row 14 onward.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F3").Value  = 1364
$ws1.Range("F7").Value  = 11700
$ws1.Range("F8").Value  = 4391
$ws1.Range("F10").Value = 40
$ws1.Range("F13").Value = 2548
$ws1.Range("F14").Value = 1095
$ws1.Range("F15").Value = 150
$ws1.Range("F16").Value = 41
$ws1.Range("F17").Value = 5105
$ws1.Range("F19").Value = 184
$ws1.Range("F20").Value = 517
$ws1.Range("F21").Value = 11346
$ws1.Range("F22").Value = 11275

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F3").Value  = 1364
$ws4.Range("F7").Value  = 11700
$ws4.Range("F8").Value  = 4391
$ws4.Range("F10").Value = 40
$ws4.Range("F13").Value = 2548
$ws4.Range("F15").Value = 1095
$ws4.Range("F16").Value = 150
$ws4.Range("F17").Value = 41
$ws4.Range("F18").Value = 5105
$ws4.Range("F20").Value = 184
$ws4.Range("F21").Value = 517
$ws4.Range("F22").Value = 11346
$ws4.Range("F23").Value = 11275
